$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin / Link / Volume(1h) columns (plain text values)
$ws.Range("E2").Value = '  -2.48%  '
$ws.Range("E3").Value = '  -3.04%  '
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  -1.24%  '
$ws.Range("E8").Value = '  -1.24%  '
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("E10").Value = '  -3.43%  '
$ws.Range("E11").Value = '  -2.18%  '
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("E13").Value = '  -2.94%  '
$ws.Range("E14").Value = '  -2.88%  '
$ws.Range("E15").Value = '  -2.04%  '
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("E17").Value = '  -1.50%  '
$ws.Range("E18").Value = '  -2.78%  '
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("E20").Value = '  -3.14%  '
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("E22").Value = '  -2.36%  '
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("E25").Value = '  -3.04%  '
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("E27").Value = '  -4.29%  '
$ws.Range("E28").Value = '  -1.89%  '
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("E30").Value = '  -11.69%  '
$ws.Range("E31").Value = '  +0.90%  '
$ws.Range("E32").Value = '  -4.47%  '
$ws.Range("E33").Value = '  -4.52%  '
$ws.Range("E34").Value = '  -4.84%  '
$ws.Range("E35").Value = '  -3.94%  '
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("E37").Value = '  -1.32%  '
$ws.Range("E38").Value = '  -2.51%  '
$ws.Range("E39").Value = '  -4.73%  '
$ws.Range("E40").Value = '  -3.55%  '
$ws.Range("E41").Value = '  -3.67%  '
$ws.Range("E42").Value = '  -9.56%  '
$ws.Range("E43").Value = '  -7.12%  '
$ws.Range("E44").Value = '  -5.48%  '
$ws.Range("E45").Value = '  +0.64%  '
$ws.Range("E46").Value = '  -1.57%  '
$ws.Range("E47").Value = '  -4.48%  '
$ws.Range("E48").Value = '  -4.35%  '
$ws.Range("E49").Value = '  -4.61%  '
$ws.Range("E50").Value = '  -3.01%  '
$ws.Range("E51").Value = '  -0.69%  '

# Update Price column values, forcing text storage (values are numeric-looking
# strings like '1.006' or multi-dot thousand-separated strings like '26.818.37')
# so they keep their original inline-string/text representation instead of being
# coerced into Excel numbers.
$ws.Range("D2").Formula = "'26.818.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Formula = "'1.776.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Formula = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Formula = "'310.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Formula = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Formula = "'0.4238"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Formula = "'0.3618"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Formula = "'0.07184"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Formula = "'0.8367"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Formula = "'20.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Formula = "'1.807.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Formula = "'5.252"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Formula = "'6.340"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Formula = "'0.06788"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Formula = "'1.006"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Formula = "'79.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Formula = "'0.000008676"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Formula = "'1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Formula = "'14.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Formula = "'26.886.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Formula = "'5.017"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Formula = "'11.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Formula = "'2.009.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Formula = "'1.920"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Formula = "'153.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Formula = "'18.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Formula = "'5.043"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Formula = "'114.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Formula = "'1.624"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Formula = "'0.08942"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Formula = "'0.7214"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Formula = "'2.845"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Formula = "'4.326"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Formula = "'1.093"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Formula = "'1.006"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Formula = "'1.074"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Formula = "'0.01892"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Formula = "'0.05082"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Formula = "'0.4917"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Formula = "'0.1607"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Formula = "'2.531"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Formula = "'6.107"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Formula = "'7.923"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Formula = "'104.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Formula = "'10.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Formula = "'0.06222"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Formula = "'0.4476"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Formula = "'1.571"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Formula = "'1.724"
$ws.Range("D51").Style = "Normal"
